$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.966.62"
$ws.Cells.Item(2, 5).Value = "  +0.20%  "
$ws.Cells.Item(3, 4).Value = "1.555.70"
$ws.Cells.Item(3, 5).Value = "  +0.31%  "
$ws.Cells.Item(4, 5).Value = "  -0.52%  "
$ws.Cells.Item(5, 4).Value = "'207.05"
$ws.Cells.Item(5, 5).Value = "  -0.07%  "
$ws.Cells.Item(6, 4).Value = "'0.490"
$ws.Cells.Item(6, 5).Value = "  +0.88%  "
$ws.Cells.Item(7, 5).Value = "  -0.54%  "
$ws.Cells.Item(8, 4).Value = "'22.10"
$ws.Cells.Item(8, 5).Value = "  +2.20%  "
$ws.Cells.Item(9, 5).Value = "  -0.21%  "
$ws.Cells.Item(10, 4).Value = "'0.0596"
$ws.Cells.Item(10, 5).Value = "  +1.51%  "
$ws.Cells.Item(11, 5).Value = "  -0.54%  "
$ws.Cells.Item(12, 4).Value = "1.776.60"
$ws.Cells.Item(12, 5).Value = "  +0.27%  "
$ws.Cells.Item(13, 4).Value = "1.552.44"
$ws.Cells.Item(13, 5).Value = "  +0.25%  "
$ws.Cells.Item(15, 5).Value = "  +1.13%  "
$ws.Cells.Item(16, 4).Value = "26.953.26"
$ws.Cells.Item(16, 5).Value = "  +0.11%  "
$ws.Cells.Item(17, 4).Value = "'61.78"
$ws.Cells.Item(17, 5).Value = "  -0.04%  "
$ws.Cells.Item(18, 2).Value = "BitcoinCash"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(18, 4).Value = "'218.02"
$ws.Cells.Item(18, 5).Value = "  +1.32%  "
$ws.Cells.Item(19, 2).Value = "ShibaInu"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(19, 4).Value = "0.0₃0707"
$ws.Cells.Item(19, 5).Value = "  +2.78%  "
$ws.Cells.Item(21, 5).Value = "  -0.57%  "
$ws.Cells.Item(22, 5).Value = "  +1.86%  "
$ws.Cells.Item(23, 5).Value = "  +0.82%  "
$ws.Cells.Item(24, 5).Value = "  -2.00%  "
$ws.Cells.Item(25, 5).Value = "  +0.40%  "
$ws.Cells.Item(26, 5).Value = "  -0.29%  "
$ws.Cells.Item(27, 5).Value = "  +1.00%  "
$ws.Cells.Item(28, 5).Value = "  +0.95%  "
$ws.Cells.Item(29, 5).Value = "  -0.51%  "
$ws.Cells.Item(30, 4).Value = "'0.0469"
$ws.Cells.Item(31, 5).Value = "  -0.31%  "
$ws.Cells.Item(32, 5).Value = "  +0.45%  "
$ws.Cells.Item(33, 4).Value = "'3.12"
$ws.Cells.Item(33, 5).Value = "  +4.09%  "
$ws.Cells.Item(34, 4).Value = "1.420.16"
$ws.Cells.Item(34, 5).Value = "  +1.04%  "
$ws.Cells.Item(35, 2).Value = "TrustWalletToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(35, 4).Value = "'1.07"
$ws.Cells.Item(35, 5).Value = "  +12.51%  "
$ws.Cells.Item(36, 2).Value = "LidoDAOToken"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(36, 4).Value = "'1.61"
$ws.Cells.Item(36, 5).Value = "  +3.13%  "
$ws.Cells.Item(37, 5).Value = "  +0.37%  "
$ws.Cells.Item(39, 4).Value = "'0.530"
$ws.Cells.Item(39, 5).Value = "  +1.85%  "
$ws.Cells.Item(40, 5).Value = "  -0.08%  "
$ws.Cells.Item(41, 5).Value = "  -0.56%  "
$ws.Cells.Item(42, 5).Value = "  +2.57%  "
$ws.Cells.Item(43, 4).Value = "'2.32"
$ws.Cells.Item(43, 5).Value = "  +2.32%  "
$ws.Cells.Item(44, 4).Value = "'0.999"
$ws.Cells.Item(44, 5).Value = "  +0.81%  "
$ws.Cells.Item(45, 4).Value = "'64.55"
$ws.Cells.Item(45, 5).Value = "  +1.52%  "
$ws.Cells.Item(46, 5).Value = "  +0.64%  "
$ws.Cells.Item(47, 4).Value = "1.690.50"
$ws.Cells.Item(47, 5).Value = "  +0.26%  "
$ws.Cells.Item(48, 4).Value = "'87.31"
$ws.Cells.Item(48, 5).Value = "  +1.29%  "
$ws.Cells.Item(49, 5).Value = "  +1.28%  "
$ws.Cells.Item(50, 5).Value = "  +0.28%  "
$ws.Cells.Item(51, 5).Value = "  +0.73%  "
